# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 20:05"

# --- Plain data refreshes (no row re-ordering) ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1440427
$ws.Cells.Item(4, 3).Value = 10079
$ws.Cells.Item(4, 4).Value = 311721
$ws.Cells.Item(4, 5).Value = 1042715
$ws.Cells.Item(4, 6).Value = 16337
$ws.Cells.Item(4, 7).Value = 794
$ws.Cells.Item(4, 8).Value = 85991

# Row 10: Francia
$ws.Cells.Item(10, 2).Value = 178060
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 59605
$ws.Cells.Item(10, 5).Value = 91030
$ws.Cells.Item(10, 6).Value = 2299
$ws.Cells.Item(10, 7).Value = 351
$ws.Cells.Item(10, 8).Value = 27425

# Row 31: Irlanda
$ws.Cells.Item(31, 2).Value = 23827
$ws.Cells.Item(31, 3).Value = 426
$ws.Cells.Item(31, 4).Value = 19470
$ws.Cells.Item(31, 5).Value = 2851
$ws.Cells.Item(31, 6).Value = 69
$ws.Cells.Item(31, 7).Value = 9
$ws.Cells.Item(31, 8).Value = 1506

# Row 36: Israel
$ws.Cells.Item(36, 2).Value = 16579
$ws.Cells.Item(36, 3).Value = 31
$ws.Cells.Item(36, 4).Value = 12521
$ws.Cells.Item(36, 5).Value = 3793
$ws.Cells.Item(36, 6).Value = 62
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 265

# Row 150: Suazilandia
$ws.Cells.Item(150, 2).Value = 187
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 51
$ws.Cells.Item(150, 5).Value = 134
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 2

# --- Countries re-sorted: Maldivas now sorts before Letonia ---
# Row 102 becomes Maldivas, row 103 becomes Letonia (values follow the swap).
$ws.Cells.Item(102, 1).Value = "Maldivas"
$ws.Cells.Item(102, 2).Value = 968
$ws.Cells.Item(102, 3).Value = 13
$ws.Cells.Item(102, 4).Value = 40
$ws.Cells.Item(102, 5).Value = 924
$ws.Cells.Item(102, 6).Value = 2
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 4

$ws.Cells.Item(103, 1).Value = "Letonia"
$ws.Cells.Item(103, 2).Value = 962
$ws.Cells.Item(103, 3).Value = 11
$ws.Cells.Item(103, 4).Value = 627
$ws.Cells.Item(103, 5).Value = 316
$ws.Cells.Item(103, 6).Value = 2
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 19

# Row 104: Sri Lanka (name unchanged, values refreshed)
$ws.Cells.Item(104, 2).Value = 925
$ws.Cells.Item(104, 3).Value = 32
$ws.Cells.Item(104, 4).Value = 445
$ws.Cells.Item(104, 5).Value = 471
$ws.Cells.Item(104, 6).Value = 1
$ws.Cells.Item(104, 7).Value = 0
$ws.Cells.Item(104, 8).Value = 9

# --- Countries re-sorted: Mali now sorts right after Costa Rica ---
# Row 112 becomes Mali, rows 113-115 shift down (Burkina Faso, Andorra, Kenia).
$ws.Cells.Item(112, 1).Value = "Mali"
$ws.Cells.Item(112, 2).Value = 779
$ws.Cells.Item(112, 3).Value = 21
$ws.Cells.Item(112, 4).Value = 436
$ws.Cells.Item(112, 5).Value = 297
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 2
$ws.Cells.Item(112, 8).Value = 46

$ws.Cells.Item(113, 1).Value = "Burkina Faso"
$ws.Cells.Item(113, 2).Value = 773
$ws.Cells.Item(113, 3).Value = 0
$ws.Cells.Item(113, 4).Value = 592
$ws.Cells.Item(113, 5).Value = 130
$ws.Cells.Item(113, 6).Value = 0
$ws.Cells.Item(113, 7).Value = 0
$ws.Cells.Item(113, 8).Value = 51

$ws.Cells.Item(114, 1).Value = "Principado de Andorra"
$ws.Cells.Item(114, 2).Value = 760
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 4).Value = 576
$ws.Cells.Item(114, 5).Value = 135
$ws.Cells.Item(114, 6).Value = 3
$ws.Cells.Item(114, 7).Value = 0
$ws.Cells.Item(114, 8).Value = 49

$ws.Cells.Item(115, 1).Value = "Kenia"
$ws.Cells.Item(115, 2).Value = 758
$ws.Cells.Item(115, 3).Value = 21
$ws.Cells.Item(115, 4).Value = 284
$ws.Cells.Item(115, 5).Value = 432
$ws.Cells.Item(115, 6).Value = 1
$ws.Cells.Item(115, 7).Value = 2
$ws.Cells.Item(115, 8).Value = 42
